$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.922.49"
$ws.Range("E2").Value = "  +0.39%  "

$ws.Range("D3").Value = "'1.874.88"
$ws.Range("E3").Value = "  +0.88%  "

$ws.Range("E4").Value = "  -0.36%  "

$ws.Range("D5").Value = "'306.54"
$ws.Range("E5").Value = "  +0.48%  "

$ws.Range("D6").Value = "'0.9991"
$ws.Range("E6").Value = "  -0.41%  "

$ws.Range("D7").Value = "'0.5157"
$ws.Range("E7").Value = "  +1.60%  "

$ws.Range("D8").Value = "'0.3738"
$ws.Range("E8").Value = "  +1.99%  "

$ws.Range("D9").Value = "'0.07191"
$ws.Range("E9").Value = "  +0.97%  "

$ws.Range("E10").Value = "  +1.77%  "

$ws.Range("D11").Value = "'20.73"
$ws.Range("E11").Value = "  +0.62%  "

$ws.Range("E12").Value = "  +0.44%  "

$ws.Range("D13").Value = "'94.93"
$ws.Range("E13").Value = "  +5.35%  "

$ws.Range("D14").Value = "'1.846.24"
$ws.Range("E14").Value = "  -0.86%  "

$ws.Range("D15").Value = "'5.256"
$ws.Range("E15").Value = "  +0.34%  "

$ws.Range("D16").Value = "'0.9993"
$ws.Range("E16").Value = "  -0.34%  "

$ws.Range("D17").Value = "'0.000008522"
$ws.Range("E17").Value = "  -0.33%  "

$ws.Range("D18").Value = "'14.26"
$ws.Range("E18").Value = "  +1.97%  "

$ws.Range("D19").Value = "'0.9989"
$ws.Range("E19").Value = "  -0.60%  "

$ws.Range("D20").Value = "'26.954.07"
$ws.Range("E20").Value = "  +0.37%  "

$ws.Range("D21").Value = "'5.035"
$ws.Range("E21").Value = "  +0.87%  "

$ws.Range("D22").Value = "'2.099.15"
$ws.Range("E22").Value = "  -0.83%  "

$ws.Range("E23").Value = "  +0.70%  "

$ws.Range("D24").Value = "'6.421"
$ws.Range("E24").Value = "  -0.28%  "

$ws.Range("D25").Value = "'146.12"
$ws.Range("E25").Value = "  +0.42%  "

$ws.Range("D26").Value = "'1.780"
$ws.Range("E26").Value = "  -2.29%  "

$ws.Range("D27").Value = "'17.99"
$ws.Range("E27").Value = "  +0.93%  "

$ws.Range("D28").Value = "'2.111"
$ws.Range("E28").Value = "  +3.31%  "

$ws.Range("D29").Value = "'114.95"
$ws.Range("E29").Value = "  +2.18%  "

$ws.Range("D30").Value = "'4.910"
$ws.Range("E30").Value = "  +5.36%  "

$ws.Range("D31").Value = "'4.775"
$ws.Range("E31").Value = "  +4.00%  "

$ws.Range("D32").Value = "'0.09186"
$ws.Range("E32").Value = "  -0.03%  "

$ws.Range("D33").Value = "'0.05028"
$ws.Range("E33").Value = "  -1.40%  "

$ws.Range("D34").Value = "'0.7527"
$ws.Range("E34").Value = "  +3.56%  "

$ws.Range("D35").Value = "'1.172"
$ws.Range("E35").Value = "  +2.56%  "

$ws.Range("D36").Value = "'2.988"
$ws.Range("E36").Value = "  -2.84%  "

$ws.Range("D37").Value = "'3.252"
$ws.Range("E37").Value = "  +1.80%  "

$ws.Range("D38").Value = "'0.5583"
$ws.Range("E38").Value = "  +5.60%  "

$ws.Range("E39").Value = "  -1.21%  "

$ws.Range("D40").Value = "'2.486"
$ws.Range("E40").Value = "  +1.30%  "

$ws.Range("E41").Value = "  -0.17%  "

$ws.Range("D42").Value = "'6.580"
$ws.Range("E42").Value = "  +2.15%  "

$ws.Range("D43").Value = "'115.67"
$ws.Range("E43").Value = "  -1.08%  "

$ws.Range("D44").Value = "'8.690"
$ws.Range("E44").Value = "  +3.61%  "

$ws.Range("D45").Value = "'0.1500"
$ws.Range("E45").Value = "  +2.34%  "

$ws.Range("D46").Value = "'0.4775"
$ws.Range("E46").Value = "  +3.69%  "

$ws.Range("D47").Value = "'0.9988"
$ws.Range("E47").Value = "  -0.33%  "

$ws.Range("D48").Value = "'10.13"
$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("D49").Value = "'1.565"
$ws.Range("E49").Value = "  +1.01%  "

$ws.Range("D50").Value = "'37.06"
$ws.Range("E50").Value = "  +0.62%  "

$ws.Range("D51").Value = "'63.41"
$ws.Range("E51").Value = "  +0.64%  "
